$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ur = $ws.UsedRange
$firstRow = $ur.Row
$firstCol = $ur.Column
$rowCount = $ur.Rows.Count
$colCount = $ur.Columns.Count
$lastRow = $firstRow + $rowCount - 1
$lastCol = $firstCol + $colCount - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -ne $null -and $v -is [string]) {
            $nv = $v.Replace("D80", "D86").Replace("D64", "D69").Replace("D51", "D55").Replace("S30", "S31")
            if ($nv -ne $v) {
                $cell.Value2 = $nv
            }
        }
    }
}
